# Rebuild flow map logic and data pipeline
# - Rename the raw/coded column headers (row 1) to human readable labels.
# - Strip the bold/bordered/centered header formatting back to the plain
#   default cell style (no bold font, no border, no special alignment).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header relabeling (A1 "state" and L1 "Year" stay the same).
$ws.Range("B1").Value = "Federal Contracts"
$ws.Range("C1").Value = "Federal Contracts (Indirect)"
$ws.Range("D1").Value = "Sub-contract Out"
$ws.Range("E1").Value = "Sub-Contract In"
$ws.Range("F1").Value = "Net Sub-Contract"
$ws.Range("G1").Value = "Employees"
$ws.Range("H1").Value = "Resident"
$ws.Range("I1").Value = "Federal Contracts per 1000 residents"
$ws.Range("J1").Value = "Federal Contracts (Indirect) per 1000 residents"
$ws.Range("K1").Value = "Net Sub-Contract per 1000 residents"

# Remove the bold font, border and centered/top alignment that used to be
# applied to the header row, returning those cells to the workbook's
# default (unstyled) appearance.
$header = $ws.Range("A1:L1")
$header.ClearFormats()
